# Adds a new "Address" column (F) to the sheet, shifting the existing
# "District" column (F) one position right to G. The address text is the
# school/place portion of each teacher's address (column B / E), with the
# trailing ", <District>." stripped and internal ", " separators removed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before F; old F (District) becomes G.
$ws.Columns("F:F").Insert()

# Header
$ws.Range("F2").Value = "Address"

# Per-row address values (rows with no address text are left blank)
$ws.Range("F3").Value = "Govt. Adarsha Vidyalaya Hosapet"
$ws.Range("F4").Value = "G P U College T B DamHospet"
$ws.Range("F5").Value = "G H S TalurSiruguppa"
$ws.Range("F6").Value = "G H P S SusheelnagarSandur"
$ws.Range("F7").Value = "Govt. Junior CollegeH B Halli"
$ws.Range("F8").Value = "Govt. High School L N HalliplotHadagali"
$ws.Range("F10").Value = "G H S K BelagalluSiruguppa"
$ws.Range("F11").Value = "Adarsha Vidyalaya H B Halli"
$ws.Range("F12").Value = "G H S KakubalHospet"
$ws.Range("F13").Value = "Govt. High School Y Bulihal"
$ws.Range("F14").Value = "G H S Upanayakana halliHadagali"
$ws.Range("F15").Value = "G H S GaddikeriH B Halli"
$ws.Range("F16").Value = "Govt. Adarsha Vidyalaya Hosapet"
$ws.Range("F17").Value = "Govt. Higher Primary School Chaganuru"
$ws.Range("F18").Value = "G H P S BenakalBallari East"
$ws.Range("F19").Value = "G H P S DanlethpuraSandur"
$ws.Range("F20").Value = "G H S HirekolachiH Hadagali"
$ws.Range("F21").Value = "G H S KoluruBallary West"
$ws.Range("F22").Value = "G G H S SirigeriSiruppa"
$ws.Range("F23").Value = "G H P S D MallapuraSandur"
$ws.Range("F24").Value = "K P G High School TambrahalliH B Halli"
$ws.Range("F25").Value = "Govt. High SchoolVattammanahalliH B Halli"
$ws.Range("F26").Value = "K M H S BachigondanahalliH B Halli"
$ws.Range("F27").Value = "H B Halli"
$ws.Range("F28").Value = "G H S (RMSA) ChilugoduH B Halli"
$ws.Range("F29").Value = "Adarsha Vidyalaya (RMSA) Siruguppa"
$ws.Range("F30").Value = "G H S AlaburuH B Halli"
$ws.Range("F31").Value = "G H S M B AyyanahalliKudligi"
$ws.Range("F32").Value = "Govt. High SchoolGadiganurHosapete"
$ws.Range("F33").Value = "G H S (RMSA) KariganurHosapete"
$ws.Range("F34").Value = "G H P S C K HalliSandur"
$ws.Range("F35").Value = "S P M Govt. High School Hosapet"
$ws.Range("F36").Value = "Govt. High School B T GuddiKudligi"
$ws.Range("F37").Value = "V K K G P U High School Section HirehadagaliHuvinahadagali"
$ws.Range("F38").Value = "Govt. High School (RMSA) VadduSandur"
$ws.Range("F39").Value = "Govt Girls High School KamalapuraHosapete"
$ws.Range("F40").Value = "G H P S SriramashettyhalliSandur"
$ws.Range("F41").Value = "G H P S KatinakambaSandur"
$ws.Range("F43").Value = "Hosapete"
$ws.Range("F44").Value = "Govt. P U College High School Section KamalapuraHosapete"
$ws.Range("F45").Value = "Govt. PU College for Girls (High Section) Hosapete"
$ws.Range("F46").Value = "STK G H SchoolMahajanadahalliHadagali"
$ws.Range("F47").Value = "Govt. High School HateholliSiruguppa"
$ws.Range("F48").Value = "Govt. Adarsha Vidyalaya Krishna NagarSandur"
$ws.Range("F49").Value = "G H S HirebannimattiHadagali"
$ws.Range("F50").Value = "G H S NandihalliH Hadagali"
$ws.Range("F51").Value = "G H S BannikalluH B Halli"
$ws.Range("F52").Value = "G H S BudanurH Hadagali"
$ws.Range("F53").Value = "Sri Verabadreshwara High School DupadahalliKudaligi"
$ws.Range("F54").Value = "G H S A S GudiHospet"
$ws.Range("F55").Value = "T M S G H S N B PuraHuvinahadagali"
$ws.Range("F56").Value = "G H P S NandihalliSandur"
$ws.Range("F57").Value = "Govt. Girls High School (EXMPL)Hosapete"
$ws.Range("F58").Value = "G H S AnkasamudraH B Halli"
$ws.Range("F59").Value = "Govt. P U College for Girls(High School Section) Hosapete"
$ws.Range("F60").Value = "B H S AmmanakereKudligi"
$ws.Range("F61").Value = "Govt. High School Kallukambha"
